# Auto-generated Word COM-interop script implementing the target diff.
$d = $word.ActiveDocument

# Op 1: replace paragraph block COM[67..70] with 15 new paragraphs
$startPara = $d.Paragraphs(67)
$endPara = $d.Paragraphs(70)
$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$blockRange.Delete()
$anchor = $d.Paragraphs(66)
$cur = $anchor
$newItems = @(
    @{style="Heading3"; text="Marketing Strategy & Client Success"},
    @{style="Normal"; text="• Rewrote services offering for multi-million dollar advertising agency data department"},
    @{style="Normal"; text="• Restructured Decision Sciences Department to scale capabilities from small-scale data analysis to comprehensive big data operations"},
    @{style="Normal"; text="• Improved project delivery timelines by 40% through introduction of version control and Agile methodologies"},
    @{style="Normal"; text="• Revealed new insights about existing customers through spatial analysis and consumer segmentation methodologies"},
    @{style="Heading3"; text="Research & Analytics Innovation"},
    @{style="Normal"; text="• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"},
    @{style="Normal"; text="• Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors"},
    @{style="Normal"; text="• Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls"},
    @{style="Normal"; text="• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"},
    @{style="Heading3"; text="Data-Driven Marketing"},
    @{style="Normal"; text="• Led multi-million dollar market research projects involving sensitive consumer data with privacy compliance"},
    @{style="Normal"; text="• Developed advanced segmentation models using demographic, psychographic, and behavioral data"},
    @{style="Normal"; text="• Created comprehensive data visualization solutions that improved clients' understanding of complex research findings"},
    @{style="Normal"; text="• Managed national polling team of five data analysts for consumer insights and market intelligence"}
)
foreach ($item in $newItems) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    $cur.Range.Text = $item.text
    $cur.Style = $item.style
}

# Op 2: replace paragraph block COM[18..65] with 40 new paragraphs
$startPara = $d.Paragraphs(18)
$endPara = $d.Paragraphs(65)
$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$blockRange.Delete()
$anchor = $d.Paragraphs(17)
$cur = $anchor
$newItems = @(
    @{style="Heading3"; text="ANALYTICS SUPERVISOR - GSD&M, Austin, TX | November 2019 – June 2020"},
    @{style="Normal"; text="Data Department Transformation and Client Strategy"},
    @{style="Normal"; text="• Transformed the small data team into a big data engineering team, going from working on small datasets on laptops to using Hadoop Clusters and Hive on AWS"},
    @{style="Normal"; text="• Rewrote the mission and offerings of the department and drafted a plan for how it would integrate with the rest of the strategy team"},
    @{style="Normal"; text="• Managed accounts for United States Air Force, Southwest Airlines/Chase and Indeed with focus on data-driven marketing insights"},
    @{style="Normal"; text="• Introduced version control and Agile methodologies to the data team, improving project delivery timelines by 40%"},
    @{style="Normal"; text="• Managed three analysts, mentoring them in advanced market research techniques and data analysis"},
    @{style="Normal"; text="• Implemented spatial analysis and consumer segmentation methodologies that revealed new insights about existing customers"},
    @{style="Heading3"; text="DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | June 2021 – May 2023"},
    @{style="Normal"; text="Consumer Insights and Market Intelligence"},
    @{style="Normal"; text="• Conceived and led implementation of comprehensive multi-tenant data warehouse integrating consumer demographic, economic, and behavioral data"},
    @{style="Normal"; text="• Overhauled the organization's survey methodology and polling operations, significantly improving data accuracy and response rates"},
    @{style="Normal"; text="• Managed and developed one of the in-house polling teams, focusing on Random Device Engagement (RDE), text message and web panel collected surveys"},
    @{style="Normal"; text="• Worked on standardizing questions, survey instruments and call methods, along with building a meta-analytical dataset for longitudinal analysis"},
    @{style="Normal"; text="• Managed a cross-functional team of eleven data engineers and analysts, establishing best practices for research methodology and data analysis"},
    @{style="Normal"; text="• Developed advanced data pipelines for machine learning applications that enhanced consumer segmentation and predictive modeling capabilities"},
    @{style="Heading3"; text="SENIOR ANALYST - Myers Research, Washington, DC | August 2012 – February 2014"},
    @{style="Normal"; text="Market Research and Consumer Insights"},
    @{style="Normal"; text="• Designed comprehensive survey instruments for specialized voting segments and niche markets"},
    @{style="Normal"; text="• Developed sophisticated analytical products and reports that delivered actionable insights to clients"},
    @{style="Normal"; text="• Co-developed RACSO web application to manage all aspects of survey operations, from instrument design to data collection and analysis"},
    @{style="Normal"; text="• Introduced geospatial techniques to enhance market segmentation capabilities, providing clients with location-based consumer insights"},
    @{style="Normal"; text="• Standardized reporting methodologies to improve clarity and impact of research findings"},
    @{style="Normal"; text="• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research"},
    @{style="Heading3"; text="RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | August 2011 – August 2012"},
    @{style="Normal"; text="Political Marketing and Campaign Strategy"},
    @{style="Normal"; text="• Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls"},
    @{style="Normal"; text="• Used FLEEM for early quantitative research in support of Senators Martin Heinrich and Elizabeth Warren"},
    @{style="Normal"; text="• Led all aspects of survey design, implementation, data analysis, and reporting for major national studies"},
    @{style="Normal"; text="• Developed new statistical methods for boundary estimation techniques, enhancing geographic market segmentation capabilities"},
    @{style="Normal"; text="• Created comprehensive data visualization solutions that improved clients' understanding of complex research findings"},
    @{style="Normal"; text="• Provided tabular and graphical reporting with plans for interactive data exploration capabilities"},
    @{style="Heading3"; text="PROGRAMMER - Lake Research Partners, Washington, DC | April 2008 – December 2008"},
    @{style="Normal"; text="Market Research and Consumer Analysis"},
    @{style="Normal"; text="• Worked on all aspects of questionnaire design, sampling, reporting and analysis for political actors in Congressional, Senate and Presidential elections"},
    @{style="Normal"; text="• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party and affiliated actors"},
    @{style="Normal"; text="• Designed questionnaires and analyzed data for complex market research studies across diverse industries"},
    @{style="Normal"; text="• Conducted statistical modeling and analysis to address multifaceted consumer behavior questions"},
    @{style="Normal"; text="• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"},
    @{style="Normal"; text="• Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings"}
)
foreach ($item in $newItems) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    $cur.Range.Text = $item.text
    $cur.Style = $item.style
}

# Op 3: replace text of paragraph 17 (COM 1-based)
$p = $d.Paragraphs(17)
$null = $p.Range.Find.Execute("• Built comprehensive competitive intelligence frameworks analyzing market trends, pricing strategies, and feature differentiation across political technology and data analytics sectors", $false, $false, $false, $false, $false, $true, 1, $false, "• Translated complex technical concepts into clear, actionable messaging that drives customer adoption and business growth", 2)

# Op 4: replace text of paragraph 16 (COM 1-based)
$p = $d.Paragraphs(16)
$null = $p.Range.Find.Execute("• Collaborated with cross-functional teams including engineering, sales, and customer success to align product strategy with market demands and customer feedback", $false, $false, $false, $false, $false, $true, 1, $false, "• Built and maintained client relationships across diverse industries, consistently delivering insights that drove strategic decision-making", 2)

# Op 5: replace text of paragraph 15 (COM 1-based)
$p = $d.Paragraphs(15)
$null = $p.Range.Find.Execute("• Conducted extensive customer research and segmentation analysis using survey methodology and behavioral data to develop targeted buyer personas and messaging frameworks", $false, $false, $false, $false, $false, $true, 1, $false, "• Developed and deployed custom research software that processed billions of consumer records for pattern analysis, fraud detection and entity resolution", 2)

# Op 6: replace text of paragraph 14 (COM 1-based)
$p = $d.Paragraphs(14)
$null = $p.Range.Find.Execute("• Created compelling product narratives and value propositions that translated complex technical capabilities into clear customer benefits, resulting in improved adoption rates and customer engagement", $false, $false, $false, $false, $false, $true, 1, $false, "• Led multi-million dollar market research projects involving sensitive consumer data, ensuring compliance with privacy regulations while delivering actionable insights", 2)

# Op 7: replace text of paragraph 13 (COM 1-based)
$p = $d.Paragraphs(13)
$null = $p.Range.Find.Execute("• Developed and executed go-to-market strategies for multiple SaaS platform launches, achieving thousands of active users and significant market penetration", $false, $false, $false, $false, $false, $true, 1, $false, "• Designed and implemented advanced segmentation models using demographic, psychographic, and behavioral data to identify high-value targets", 2)

# Op 8: replace text of paragraph 12 (COM 1-based)
$p = $d.Paragraphs(12)
$null = $p.Range.Find.Execute("• Led comprehensive market intelligence and competitive analysis projects for B2B technology platforms, delivering actionable insights that shaped product positioning and messaging strategies", $false, $false, $false, $false, $false, $true, 1, $false, "• Conducted comprehensive quantitative and qualitative research studies for political candidates and major organizations, providing actionable consumer insights and market intelligence", 2)

# Op 9: replace text of paragraph 11 (COM 1-based)
$p = $d.Paragraphs(11)
$null = $p.Range.Find.Execute("Market Research, Product Strategy & Go-to-Market Leadership", $false, $false, $false, $false, $false, $true, 1, $false, "Marketing Strategy and Data-Driven Insights", 2)

# Op 10: replace text of paragraph 10 (COM 1-based)
$p = $d.Paragraphs(10)
$null = $p.Range.Find.Execute("PARTNER - Siege Analytics, Washington, DC | January 2014 – Present", $false, $false, $false, $false, $false, $true, 1, $false, "PARTNER - Siege Analytics, Washington, DC | 2005 – Present", 2)

# Op 11: replace text of paragraph 8 (COM 1-based)
$p = $d.Paragraphs(8)
$null = $p.Range.Find.Execute("Communication & Technology: Strategic Messaging & Narrative Development • Stakeholder Communication & Executive Briefings • Content Creation: Case Studies, Battle Cards, Playbooks • B2B SaaS Platform Experience & Technical Acumen • CRM/Marketing Automation (Salesforce, HubSpot) • Data Visualization (Tableau, PowerBI, D3.js) • AI/ML Tools Integration & Marketing Technology Stack", $false, $false, $false, $false, $false, $true, 1, $false, "Communication & Technology: Strategic Messaging & Narrative Development • Technical Concept Translation for Business Audiences • Stakeholder Communication & Presentation Skills • Data Visualization & Reporting (Tableau, PowerBI, d3.js) • Marketing Technology Stack Integration • Content Strategy & Thought Leadership • Client Relationship Management & Business Development", 2)

# Op 12: replace text of paragraph 4 (COM 1-based)
$p = $d.Paragraphs(4)
$null = $p.Range.Find.Execute("Results-driven Product Marketing professional with 21 years of experience translating complex data insights into compelling market strategies and customer narratives. Expert in market intelligence, competitive analysis, and data-driven positioning with proven success leading cross-functional teams and launching B2B SaaS platforms used by thousands of users. Deep expertise in survey methodology, customer segmentation, and go-to-market strategy development. Skilled at turning complex technical concepts into clear, actionable messaging that drives customer adoption and business growth across political, technology, and consulting sectors.", $false, $false, $false, $false, $false, $true, 1, $false, "Results-driven Marketing & Data Analytics Professional with 21 years of experience translating complex data insights into compelling market strategies and customer narratives. Expert in market intelligence, competitive analysis, and data-driven positioning with proven success leading cross-functional teams and launching B2B SaaS platforms used by thousands of users. Deep expertise in survey methodology, customer segmentation, and go-to-market strategy development. Skilled at turning complex technical concepts into clear, actionable messaging that drives customer adoption and business growth across political, technology, and consulting sectors.", 2)

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
